$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsTests  = $wb.Worksheets.Item("Test Cases")

# --- Config sheet: update BaseURL value and turn it into a hyperlink ---
$wsConfig.Range("B5").Value = "http://testingmasters.com/hrm/symfony/web/index.php/auth/login"
$wsConfig.Hyperlinks.Add($wsConfig.Range("B5"), "http://testingmasters.com/hrm/symfony/web/index.php/auth/login")

# --- Test Cases sheet: update the Execute expression and add a new test row ---
$wsTests.Range("D2").Value = "Groups=Orange"

$wsTests.Range("A13").Value = "201"
$wsTests.Range("B13").Value = "1"
$wsTests.Range("C13").Value = "Orange"

# --- Selection / active sheet bookkeeping (matches the saved view state) ---
$wsConfig.Select()
$wsConfig.Range("B11").Select()

$wsTests.Select()
$wsTests.Range("D3").Select()

Write-Host "Updated TestData for Orange Hrm"
